# Add the "Coin Change" problem as a new row at the bottom of the "DW" sheet.
# The new row reuses the same look as the existing "Easy" rows (green
# font/fill with a thin grey border) but, since it sits at the very bottom
# of the table, only keeps the left/right edges of that border (no top or
# bottom rule). Column D keeps the sheet's plain default look, matching the
# rest of the "Link" cells that don't carry the colored highlight.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DW")

$row = 11

$ws.Cells.Item($row, 1).Value = "Coin Change"
$ws.Cells.Item($row, 2).Value = "Dynamic programming"
$ws.Cells.Item($row, 3).Value = "Medium"
$ws.Cells.Item($row, 4).Value = "https://leetcode.com/problems/coin-change/"
$ws.Cells.Item($row, 5).Value = "Dynamic programming bottm up approach start by calculating dp to complete 1 ruppe then 2 ruppes then so onn "

# Build the final look (green font/fill + thin left/right-only border) once
# on an out-of-the-way scratch cell, so the working sheet only ever receives
# the finished format in a single paste per target cell.
$scratch = $ws.Range("Z100")
$ws.Range("A5").Copy()
$scratch.PasteSpecial(-4122)
$scratch.Borders.Item(9).LineStyle = -4142
$scratch.Borders.Item(8).LineStyle = -4142

$scratch.Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("E11").PasteSpecial(-4122)

$scratch.Clear()

# Column D (the link) keeps the sheet's default column formatting - no
# explicit style needed.

$ws.Activate()
$ws.Range("E11").Select()
